# Update Name of Algo
# Apply updated KNN imputation result values to columns C and D for the
# affected rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.862
$ws.Range("C9").Value = -11.002
$ws.Range("D9").Value = -7.356
$ws.Range("D11").Value = -7.323
$ws.Range("C13").Value = -13.529
$ws.Range("C16").Value = -12.888
$ws.Range("D16").Value = -8.1
$ws.Range("C18").Value = -11.888
$ws.Range("C20").Value = -12.771
$ws.Range("D23").Value = -8.085999999999999
$ws.Range("D24").Value = -7.401999999999998
$ws.Range("C26").Value = -11.985
$ws.Range("D26").Value = -7.530999999999999
$ws.Range("C27").Value = -12.486
$ws.Range("C29").Value = -12.288
$ws.Range("D34").Value = -7.516999999999999
$ws.Range("C35").Value = -12.428
$ws.Range("D35").Value = -8.088999999999999
$ws.Range("C36").Value = -12.401
$ws.Range("D44").Value = -7.763
$ws.Range("C45").Value = -12.835
$ws.Range("D48").Value = -7.904000000000001
$ws.Range("D49").Value = -8.097
$ws.Range("D52").Value = -7.896000000000001
$ws.Range("C55").Value = -13.585
$ws.Range("C57").Value = -13.643
$ws.Range("D66").Value = -7.577
$ws.Range("D67").Value = -7.789
$ws.Range("C69").Value = -11.105
$ws.Range("D73").Value = -8.19
$ws.Range("C76").Value = -12.994
$ws.Range("C78").Value = -12.296
$ws.Range("D78").Value = -7.486
$ws.Range("D80").Value = -8.282999999999998
$ws.Range("C82").Value = -11.9
$ws.Range("C83").Value = -13.268
$ws.Range("D91").Value = -7.376
$ws.Range("C93").Value = -10.818
$ws.Range("C97").Value = -12.008
$ws.Range("D97").Value = -7.597
$ws.Range("D99").Value = -7.968000000000001
$ws.Range("D104").Value = -7.790999999999999
